# The diff shows every existing record from row 101 downward shifting
# down by one row (old row 101 -> new row 102, ..., old row 142 -> new
# row 143), while a brand-new record is inserted at row 101. This is a
# classic "insert a row above" edit (weekly data refresh adding the
# newest price observation at the top of the Arveja Verde records),
# so we do exactly that via the Excel object model instead of rewriting
# every row by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row at row 101; Excel shifts rows 101:142 down to
# 102:143 and extends the used range accordingly.
$ws.Rows("101:101").Insert()

# Populate the newly inserted row 101 with the new record's data.
$ws.Range("A101").Value = 4
$ws.Range("B101").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C101").Value = "Los Lagos"
$ws.Range("D101").Value = 44845
$ws.Range("E101").Value = 10
$ws.Range("F101").Value = 100112022
$ws.Range("G101").Value = "Arveja Verde"
$ws.Range("H101").Value = "Perfection"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 70
$ws.Range("K101").Value = 34000
$ws.Range("L101").Value = 34000
$ws.Range("M101").Value = 34000
$ws.Range("N101").Value = '$/malla 25 kilos'
$ws.Range("O101").Value = "Provincia de Huasco"
$ws.Range("P101").Value = 1360
$ws.Range("Q101").Value = 25
$ws.Range("R101").Value = "Hortaliza"
